$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record for "Haba" (Primera, Provincia de Limarí) needs to be
# inserted as row 199, pushing the existing rows 199-254 down to 200-255.
$ws.Rows(199).Insert()

$ws.Range("A199").Value = 9
$ws.Range("B199").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C199").Value = "Metropolitana"
$ws.Range("D199").Value = 44785
$ws.Range("E199").Value = 13
$ws.Range("F199").Value = 100112026
$ws.Range("G199").Value = "Haba"
$ws.Range("H199").Value = "Sin especificar"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 52
$ws.Range("K199").Value = 13000
$ws.Range("L199").Value = 14000
$ws.Range("M199").Value = 13500
$ws.Range("N199").Value = "`$/saco 25 kilos"
$ws.Range("O199").Value = "Provincia de Limarí"
$ws.Range("P199").Value = 540
$ws.Range("Q199").Value = 25
$ws.Range("R199").Value = "Hortaliza"
